$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 2D training schedule values (rows 2-6, columns A-J).
# Column order: trialTrain, x_fixStart, y_fixStart, x_corrSteps, y_corrSteps,
#               x_nrSteps, y_nrSteps, alienID, praclen, version
$data = @(
    @(1, 8, 7, 6, 3, -2, -4, 32, 5, "train_dim2_1"),
    @(2, 7, 7, 6, 2, -1, -5, 21, 5, "train_dim2_1"),
    @(3, 6, 8, 1, 7, -5, -1, 65, 5, "train_dim2_1"),
    @(4, 9, 9, 5, 7, -4, -2, 54, 5, "train_dim2_1"),
    @(5, 7, 5, 4, 2, -3, -3, 43, 5, "train_dim2_1")
)

$row = 2
foreach ($rowValues in $data) {
    $col = 1
    foreach ($val in $rowValues) {
        $ws.Cells.Item($row, $col).Value = $val
        $col++
    }
    $row++
}

# Select I1, matching the saved view state in the workbook.
$ws.Range("I1").Select()

$wb.Save()
